$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B updates (RandomForest imputed values changed)
$ws.Range("B3").Value = 5.991199999999991
$ws.Range("B4").Value = 4.920200000000002
$ws.Range("B7").Value = 6.579800000000001
$ws.Range("B8").Value = 5.379199999999997
$ws.Range("B12").Value = 5.822099999999999
$ws.Range("B14").Value = 8.657300000000006
$ws.Range("B22").Value = 5.753100000000003

# Column A updates
$ws.Range("A11").Value = -21.93510000000002
$ws.Range("A12").Value = -22.66790000000001
$ws.Range("A15").Value = -21.36810000000002
